$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D3 value (Kadastro Müdür Yrd. Mevcut_Sayisi) from 1 to 2
$ws.Range("D3").Value = 2

# Update selection to K24
$ws.Range("K24").Select()


